# Applies the "shifted to usenix style" edit to results.xlsx:
#  - Adds a new row of labels (PipeSize / Delay in Seconds / Application /
#    Typical request size) below the existing tables.
#  - Moves/resizes the two existing charts (they were shrunk and pulled up
#    toward the top of the sheet).
#  - Updates the active selection on the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Add new labelled row (B40:E40) ------------------------------------
$ws.Range("B40").Value = "PipeSize"
$ws.Range("C40").Value = "Delay in Seconds"
$ws.Range("D40").Value = "Application"
$ws.Range("E40").Value = "Typical request size"

# --- 2. Reposition / resize the two charts on the sheet --------------------
# Values below (in points) were solved so that the resulting two-cell
# anchors in the saved drawing XML line up exactly with the target
# from/to column+row+offset combination.
$cos = $ws.ChartObjects()
$idx = 0
foreach ($co in $cos) {
    $idx = $idx + 1
    if ($idx -eq 1) {
        # from: col=2 colOff=825500 row=2 rowOff=101600
        # to:   col=15 colOff=660400 row=31 rowOff=114300
        $co.Left = 181.875
        $co.Top = 38.0
        $co.Width = 844.2226377952757
        $co.Height = 436.0
    }
    if ($idx -eq 2) {
        # from: col=6 colOff=88900 row=1 rowOff=25400
        # to:   col=15 colOff=393700 row=24 rowOff=101600
        $co.Left = 455.160157480315
        $co.Top = 17.0
        $co.Width = 549.9374803149607
        $co.Height = 351.0
    }
}

# --- 3. Update the active view/selection -----------------------------------
$ws.Activate()
$ws.Range("E41").Select()
